# sua chiet khau cua sale phu va update chien luoc chay tinh luong theo gio
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Đơn phụ phẫu 1" — add a new service row (row 7) for a new order,
#    pushing the previous "Tổng" (total) row down to row 8 and bumping
#    its aggregated totals.
# ---------------------------------------------------------------------
$wsPhu = $wb.Worksheets.Item("Đơn phụ phẫu 1")

# Insert a fresh row above the current "Tổng" row (row 7); this shifts
# the total row down to row 8, same as inserting a new order above it.
$wsPhu.Rows.Item(7).Insert()

$wsPhu.Range("A7").Value = "HD-LUXURY"
$wsPhu.Range("B7").Value = 574

# Column C holds dates formatted as plain text (e.g. "07-01-2024"), so
# force a text format before typing the value, otherwise Excel's
# auto-detection would silently convert it into a real date serial.
$wsPhu.Range("C7").NumberFormat = "@"
$wsPhu.Range("C7").Value = "07-21-2024"
$wsPhu.Range("C7").ClearFormats()

$wsPhu.Range("D7").Value = "LONG XUYÊN"
$wsPhu.Range("E7").Value = "Trần thị mỹ duyên"
$wsPhu.Range("F7").Value = "Cá nhân"
$wsPhu.Range("G7").Value = "Cắt mí"
$wsPhu.Range("H7").Value = "Đào Vương Anh"
$wsPhu.Range("I7").Value = 50000

# "Tổng" row, now on row 8: bump the order count and the total amount.
$wsPhu.Range("B8").Value = 6
$wsPhu.Range("I8").Value = 250000

# ---------------------------------------------------------------------
# 2) "Lương" — update hourly-based salary run: more công (work hours/
#    days) and a higher phụ cấp (allowance), rippling into the LONG
#    XUYÊN base salary and the total salary figures.
# ---------------------------------------------------------------------
$wsLuong = $wb.Worksheets.Item("Lương")

$wsLuong.Range("B2").Value = 20.5                  # Ngày công
$wsLuong.Range("B3").Value = 717500                 # Phụ cấp
$wsLuong.Range("B12").Value = 2196428.571428571     # Lương cơ bản tại LONG XUYÊN
$wsLuong.Range("B17").Value = 250000                 # Công phụ phẫu 1 tại LONG XUYÊN
$wsLuong.Range("B29").Value = 3163928.571428571     # Tổng lương tại LONG XUYÊN
$wsLuong.Range("B31").Value = 3163928.571428571     # Tổng lương
